$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.314.83"
$ws.Range("E2").Value = "  +3.85%  "
$ws.Range("D3").Value = "1.592.23"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "213.99"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "24.05"
$ws.Range("E8").Value = "  +8.48%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("D10").Value = "0.0600"
$ws.Range("E10").Value = "  +0.83%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "1.821.66"
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "1.583.43"
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("D16").Value = "28.335.56"
$ws.Range("E16").Value = "  +4.06%  "
$ws.Range("D17").Value = "63.10"
$ws.Range("E17").Value = "  +2.02%  "
$ws.Range("D18").Value = "227.46"
$ws.Range("E18").Value = "  +4.35%  "
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").Value = "7.47"
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("E24").Value = "  +0.25%  "
$ws.Range("D25").Value = "151.59"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").Value = "15.19"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "0.0474"
$ws.Range("E31").Value = "  +0.80%  "
$ws.Range("D32").Value = "3.23"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "1.399.14"
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("D35").Value = "1.59"
$ws.Range("E35").Value = "  -1.88%  "
$ws.Range("D36").Value = "1.04"
$ws.Range("E36").Value = "  -6.20%  "
$ws.Range("E37").Value = "  +0.41%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").Value = "2.56"
$ws.Range("E39").Value = "  +9.10%  "
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "1.89"
$ws.Range("E43").Value = "  +7.59%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "5.66"
$ws.Range("E44").Value = "  -3.64%  "
$ws.Range("E45").Value = "  +0.17%  "
$ws.Range("D46").Value = "64.26"
$ws.Range("E46").Value = "  -0.31%  "
$ws.Range("D47").Value = "1.732.17"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "87.54"
$ws.Range("E48").Value = "  +1.93%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("E51").Value = "  +0.45%  "
